$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(7921, 610.1947, 0, 488.15576),
    @(14073, 608.8172, 0, 487.05376),
    @(6605, 608.8172, 0, 487.05376),
    @(2780, 606.3536, 0, 485.08288),
    @(1279, 524.57925, 0, 419.6634),
    @(6059, 512.1137, 0, 409.69096),
    @(15097, 511.92305, 0, 409.53844),
    @(12717, 507.4982, 0, 405.99856),
    @(5571, 505.466975, 0, 404.37358),
    @(16328, 503.2411, 0, 402.59288),
    @(1791, 502.386775, 0, 401.90942),
    @(5799, 499.1681, 0, 399.33448),
    @(15557, 498.798, 0, 399.0384),
    @(9681, 498.47285, 0, 398.77828),
    @(6374, 497.546725, 0, 398.03738),
    @(20033, 497.1067, 0, 397.68536),
    @(805, 496.7993750000001, 0, 397.4395000000001),
    @(2733, 494.807625, 0, 395.8461),
    @(3092, 494.807625, 0, 395.8461),
    @(8290, 494.807625, 0, 395.8461),
    @(8903, 494.807625, 0, 395.8461),
    @(8980, 493.585275, 0, 394.8682200000001),
    @(2787, 493.339975, 0, 394.67198),
    @(8528, 493.339975, 0, 394.67198),
    @(11624, 493.339975, 0, 394.67198),
    @(21454, 493.339975, 0, 394.67198),
    @(2889, 493.339975, 0, 394.67198),
    @(1740, 493.339975, 0, 394.67198),
    @(4376, 493.339975, 0, 394.67198),
    @(12837, 493.339975, 0, 394.67198),
    @(14461, 493.339975, 0, 394.67198),
    @(247, 493.339975, 0, 394.67198),
    @(10079, 493.339975, 0, 394.67198),
    @(5825, 493.339975, 0, 394.67198),
    @(4922, 493.339975, 0, 394.67198),
    @(12983, 493.339975, 0, 394.67198),
    @(6097, 483.473575, 0, 386.77886),
    @(10333, 483.119375, 0, 386.4955),
    @(1920, 480.058375, 0, 384.0467),
    @(18516, 480.058375, 0, 384.0467),
    @(13597, 480.058375, 0, 384.0467),
    @(5836, 479.588875, 0, 383.6711),
    @(7750, 478.9514, 0, 383.16112),
    @(3547, 478.728175, 0, 382.98254),
    @(13538, 477.2023, 0, 381.7618400000001),
    @(9518, 476.428975, 0, 381.14318),
    @(5739, 475.668575, 0, 380.53486),
    @(2138, 475.214225, 0, 380.17138),
    @(11720, 475.212475, 0, 380.16998),
    @(14221, 474.736775, 0, 379.78942)
)

$r = 2
foreach ($row in $rows) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
